$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new task row entry: "Requirements Doc" task label in A11
$ws.Range("A11").Value = "Requirements Doc"

# Record hours logged for two tasks this week
$ws.Range("E9").Value = 1
$ws.Range("E11").Value = 2

# Move the active selection to F11 (where the user left off editing)
$ws.Range("F11").Select()
